$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.120.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.319.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.23%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'302.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.10%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'99.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.84%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.29%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.62%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +4.23%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.58%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.18%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'17.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.90%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.21%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.679.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.317.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.10%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -3.12%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'43.022.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.64%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.08%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.46%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'240.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.77%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -0.41%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.09%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'25.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.28%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'168.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.86%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.08%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.82%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'33.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.36%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +6.55%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.82%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.07%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'18.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.94%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -1.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0696"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.47%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.57%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.50%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.93%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.44%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.996.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.46%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.07%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.33%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -10.86%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'17.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.30%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.89%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'76.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +9.08%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'54.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.75%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.545.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.11%  "
$ws.Range("E51").Style = "Normal"
